$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.539.83"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "3.456.64"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.453.30"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").Value = "4.039.84"
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000195"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.07%  "
$ws.Range("D17").Value = "64.534.25"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "3.429.44"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.547"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +19.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.46%  "
$ws.Range("E31").Value = "  +9.08%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.47%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0779"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.24%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "2.911.98"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0321"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.769"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.57%  "
$ws.Range("E48").Value = "  +3.64%  "
$ws.Range("E49").Value = "  +15.60%  "
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.06%  "
